$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = '@'
    $r.Value = $val
    $r.Style = 'Normal'
}

$ws.Range('D2').Value = '27.034.24'
$ws.Range('E2').Value = '  +0.28%  '

$ws.Range('D3').Value = '1.562.46'
$ws.Range('E3').Value = '  +0.41%  '

$ws.Range('E4').Value = '  +0.35%  '

Set-TextValue 'D5' '208.47'
$ws.Range('E5').Value = '  +0.64%  '

$ws.Range('E6').Value = '  +0.33%  '

Set-TextValue 'D8' '22.03'
$ws.Range('E8').Value = '  -0.35%  '

$ws.Range('E9').Value = '  +0.55%  '

$ws.Range('E10').Value = '  +1.68%  '

$ws.Range('E11').Value = '  -0.29%  '

$ws.Range('D12').Value = '1.785.19'
$ws.Range('E12').Value = '  +0.37%  '

$ws.Range('D13').Value = '1.533.26'
$ws.Range('E13').Value = '  -1.51%  '

$ws.Range('E14').Value = '  -0.12%  '

Set-TextValue 'D15' '0.520'
$ws.Range('E15').Value = '  +0.14%  '

$ws.Range('D16').Value = '27.043.14'
$ws.Range('E16').Value = '  +0.28%  '

Set-TextValue 'D17' '61.92'
$ws.Range('E17').Value = '  +0.27%  '

$ws.Range('D18').Value = '0.0₃0708'
$ws.Range('E18').Value = '  +1.48%  '

Set-TextValue 'D19' '216.25'
$ws.Range('E19').Value = '  -0.90%  '

Set-TextValue 'D20' '7.37'
$ws.Range('E20').Value = '  +0.87%  '

$ws.Range('E21').Value = '  +0.42%  '

$ws.Range('E22').Value = '  +1.88%  '

Set-TextValue 'D23' '9.21'
$ws.Range('E23').Value = '  -0.38%  '

$ws.Range('E24').Value = '  -0.22%  '

Set-TextValue 'D25' '153.44'
$ws.Range('E25').Value = '  -0.41%  '

Set-TextValue 'D26' '6.59'
$ws.Range('E26').Value = '  -0.91%  '

Set-TextValue 'D27' '15.07'
$ws.Range('E27').Value = '  +0.73%  '

$ws.Range('E28').Value = '  +1.49%  '

$ws.Range('E29').Value = '  +0.40%  '

$ws.Range('E30').Value = '  +0.85%  '

$ws.Range('E31').Value = '  +3.09%  '

$ws.Range('E32').Value = '  -0.26%  '

$ws.Range('E33').Value = '  +2.76%  '

$ws.Range('D34').Value = '1.427.91'
$ws.Range('E34').Value = '  +0.39%  '

$ws.Range('E35').Value = '  +1.33%  '

$ws.Range('E36').Value = '  +7.84%  '

Set-TextValue 'D37' '2.34'
$ws.Range('E37').Value = '  +2.41%  '

$ws.Range('E38').Value = '  +0.68%  '

Set-TextValue 'D39' '0.532'
$ws.Range('E39').Value = '  +2.13%  '

Set-TextValue 'D40' '5.89'
$ws.Range('E40').Value = '  +2.01%  '

$ws.Range('E41').Value = '  -0.52%  '

$ws.Range('E42').Value = '  +0.46%  '

$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D43' '1.00'
$ws.Range('E43').Value = '  +1.53%  '

$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D44' '2.31'
$ws.Range('E44').Value = '  -0.59%  '

Set-TextValue 'D45' '64.65'
$ws.Range('E45').Value = '  +0.38%  '

$ws.Range('E46').Value = '  -0.33%  '

$ws.Range('D47').Value = '1.701.42'

Set-TextValue 'D48' '87.04'
$ws.Range('E48').Value = '  -1.20%  '

$ws.Range('E49').Value = '  +4.63%  '

$ws.Range('E50').Value = '  -0.34%  '

$ws.Range('E51').Value = '  +0.45%  '
